$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.712.79"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.587.67"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  +0.84%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.10"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.21"
$ws.Range("E8").Value = "  -4.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.252"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0589"
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0868"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "1.812.80"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").Value = "1.573.70"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.88"
$ws.Range("E14").Value = "  -2.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.531"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").Value = "27.681.33"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.32"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "220.49"
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0693"
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("E22").Value = "  -3.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.60"
$ws.Range("E23").Value = "  -3.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.11"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.10"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("E29").Value = "  -3.30%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0468"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("D33").Value = "1.367.55"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.94"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.981"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.18"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.57"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.75"
$ws.Range("E45").Value = "  -4.68%  "
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").Value = "1.723.34"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.40"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  +9.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0498"
$ws.Range("E51").Value = "  -1.01%  "
